$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (rows 710 .. 726), matching the diff.
# G=$null means the "Localisation douleur" cell stays empty (uses the
# "empty" style template), otherwise it holds the given text.
$rows = @(
    @{R=710; B="Yoann Martelat";   C=70; D=5; E=4; F=6; G="Genou";      H=3},
    @{R=711; B="Kamal Bafounta";   C=70; D=6; E=6; F=0; G=$null;        H=7},
    @{R=712; B="Naim Ighbane";     C=70; D=5; E=6; F=0; G=$null;        H=8},
    @{R=713; B="Omar Benyounes";   C=70; D=5; E=2; F=6; G="Quadri ";    H=3},
    @{R=714; B="Maé Clavel";       C=70; D=5; E=5; F=6; G="Ischio";     H=7},
    @{R=715; B="Karim Belmahi";    C=70; D=5; E=7; F=0; G=$null;        H=10},
    @{R=716; B="Jeremie Laurent";  C=70; D=7; E=6; F=0; G=$null;        H=7},
    @{R=717; B="Yoan Zouma";       C=70; D=6; E=5; F=1; G="Dos";        H=6},
    @{R=718; B="Levy Ndoutoume";   C=70; D=7; E=7; F=0; G=$null;        H=6},
    @{R=719; B="Hedi Nasri";       C=70; D=6; E=5; F=2; G="Hanche";     H=7},
    @{R=720; B="Amine Taiar";      C=70; D=4; E=6; F=5; G="Courbature"; H=8},
    @{R=721; B="Ilan Ihaddadene";  C=70; D=6; E=5; F=0; G=$null;        H=2},
    @{R=722; B="Emmanuel Valey";   C=70; D=6; E=5; F=6; G="Cheville";   H=6},
    @{R=723; B="Karahali Souaré";  C=70; D=6; E=6; F=6; G="Cheville";   H=6},
    @{R=724; B="Theo Owono";       C=70; D=4; E=3; F=7; G="Genou";      H=9},
    @{R=725; B="Sofiane Belle";    C=70; D=6; E=4; F=2; G="Genou";      H=8},
    @{R=726; B="Romain Thunet";    C=70; D=5; E=6; F=0; G=$null;        H=3}
)

foreach ($row in $rows) {
    $r = $row.R
    if ($row.G -eq $null) {
        $tmpl = 704   # template row whose "Localisation douleur" cell is empty (style 2)
    } else {
        $tmpl = 709   # template row whose "Localisation douleur" cell has text (style 1)
    }
    $ws.Range("A${tmpl}:H${tmpl}").Copy()
    $ws.Range("A${r}:H${r}").PasteSpecial(-4122)

    $ws.Range("A${r}").Value = 46035
    $ws.Range("B${r}").Value = $row.B
    $ws.Range("C${r}").Value = $row.C
    $ws.Range("D${r}").Value = $row.D
    $ws.Range("E${r}").Value = $row.E
    $ws.Range("F${r}").Value = $row.F
    if ($row.G -ne $null) {
        $ws.Range("G${r}").Value = $row.G
    }
    $ws.Range("H${r}").Value = $row.H
}

# Extend the "Charge" shared formula down through the new rows.
$ws.Range("I710:I726").Formula = "=C710*D710"

$excel.Calculate()

# Match the author's final scroll position / selection.
$excel.ActiveWindow.ScrollRow = 695
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K721").Select()
